# Apply updated calibration estimates (std moving average) to the workbook.
$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.8325996345274745
$wsLinear.Range("B3").Value = -2.160574131199933
$wsLinear.Range("B4").Value = 5250.041118320721
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.1531552292575467, -0.0029520789047667085, -0.041004015630244806, -0.03344960914253515, -0.014257100032814472, 0.17240913337606406, 0.2894059632974835, 0.1879187042854853, 0.014182847186769832, -0.04921939038438452, -0.03016981956031537, -0.013501289182651316, 0.17623126613779153, 0.34956871190031946, 0.13852488122985424, -0.025521373524542923, -0.04776481075920933, -0.041778119301518354, -0.03172153865535137]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -1.456256747908645
$wsNonLinear.Range("B5").Value = -4.544239728252996
$wsNonLinear.Range("B6").Value = 5477.099820234803
$wsNonLinear.Range("B7").Value = -3.232994079944983
$wsNonLinear.Range("B8").Value = 3.217338853904733
$wsNonLinear.Range("B9").Value = 5012.011607033361
$wsNonLinear.Range("B10").Value = "[1.0, 0.1520341384828875, -0.0022602579685163603, -0.03961529943471662, -0.03120970704849033, -0.013513984579079864, 0.1716777250660928, 0.2866571277883445, 0.1873515835670389, 0.014436157335517803, -0.047821122922222827, -0.028514731814409278, -0.012892815693745004, 0.17599886519482436, 0.3475954608787413, 0.13757785876007295, -0.024129560676037366, -0.046026283603706676, -0.040212026850117204, -0.030529404596876853]"
